# Economic Dashboard update - 2026-01-13
# Rolls each indicator's Present/Lag1-4 window forward by one observation
# and bumps the "Latest Period" date for the rows that got new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CPI block (rows 18-21): date moves from 45962 -> 45992 and picks up
#     the highlighted "updated" fill (style 48 -> 49). Copy the fill/format
#     from a cell that already carries that style (N29) so the resulting
#     style index matches exactly, then overwrite with the new date.
$ws.Range("N29").Copy()
$ws.Range("N18:N21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("N18").Value = 45992
$ws.Range("N19").Value = 45992
$ws.Range("N20").Value = 45992
$ws.Range("N21").Value = 45992

# Row 18 (CPI, CPIAUCSL, M/M % Delta)
$ws.Range("Q18").Value = 0.003073552984176775
$ws.Range("R18").Value = ""
$ws.Range("S18").Value = ""
$ws.Range("T18").Value = 0.00310486015759337
$ws.Range("U18").Value = 0.003824519141221616

# Row 19 (CPIAUCSL, Y/Y % Delta)
$ws.Range("Q19").Value = 0.02653312468710926
$ws.Range("R19").Value = 0.0271196938527219
$ws.Range("S19").Value = ""
$ws.Range("T19").Value = 0.03022699626172379
$ws.Range("U19").Value = 0.02939219624933549

# Row 20 (Core CPI, CPILFESL, M/M % Delta)
$ws.Range("Q20").Value = 0.00239225778389951
$ws.Range("R20").Value = ""
$ws.Range("S20").Value = ""
$ws.Range("T20").Value = 0.002271121582325675
$ws.Range("U20").Value = 0.003459544325982167

# Row 21 (CPILFESL, Y/Y % Delta)
$ws.Range("Q21").Value = 0.02648965653766215
$ws.Range("R21").Value = 0.02618878615332623
$ws.Range("S21").Value = ""
$ws.Range("T21").Value = 0.03025542724453378
$ws.Range("U21").Value = 0.03112190821006822

# --- Rates block (rows 29-30): daily series, date moves 46031 -> 46034
$ws.Range("N29").Value = 46034
$ws.Range("Q29").Value = 2.22
$ws.Range("R29").Value = 2.24
$ws.Range("S29").Value = 2.23
$ws.Range("T29").Value = 2.24
$ws.Range("U29").Value = 2.24

$ws.Range("N30").Value = 46034
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.27
$ws.Range("U30").Value = 2.27

# --- Rates block (rows 47-52): daily series, date moves 46030 -> 46031
$ws.Range("N47").Value = 46031

$ws.Range("N48").Value = 46031
$ws.Range("Q48").Value = 3.54
$ws.Range("R48").Value = 3.49
$ws.Range("S48").Value = 3.47
$ws.Range("T48").Value = 3.47
$ws.Range("U48").Value = 3.46

$ws.Range("N49").Value = 46031
$ws.Range("Q49").Value = 3.75
$ws.Range("R49").Value = 3.74
$ws.Range("S49").Value = 3.7
$ws.Range("T49").Value = 3.72
$ws.Range("U49").Value = 3.71

$ws.Range("N50").Value = 46031
$ws.Range("Q50").Value = 4.18
$ws.Range("R50").Value = 4.19
$ws.Range("S50").Value = 4.15
$ws.Range("T50").Value = 4.18
$ws.Range("U50").Value = 4.17

$ws.Range("N52").Value = 46031
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.92
$ws.Range("S52").Value = 5.88
$ws.Range("T52").Value = 5.92
$ws.Range("U52").Value = 5.92

Write-Host "Dashboard update applied"
